$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A2 from "42563690@mynwu.ac.za" to "jacob.smith@example.com"
$ws.Range("A2").Value = "jacob.smith@example.com"

# Fix up the existing hyperlink on A2 to point at the correct mailto address
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address($false, $false) -eq "A2") {
        $hl.Address = "mailto:jacob.smith@example.com"
        $hl.TextToDisplay = "jacob.smith@example.com"
    }
}

# Restore the active selection to A2
$ws.Range("A2").Select()
